$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 99, shifting existing rows 99:202 down to 100:203.
$ws.Rows(99).Insert()

# Populate the newly inserted row 99 with the new weekly record.
$ws.Cells.Item(99, 1).Value = 8
$ws.Cells.Item(99, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(99, 3).Value = "Coquimbo"
$ws.Cells.Item(99, 4).Value = 45128
$ws.Cells.Item(99, 5).Value = 4
$ws.Cells.Item(99, 6).Value = 100112052
$ws.Cells.Item(99, 7).Value = "Albahaca"
$ws.Cells.Item(99, 8).Value = "Sin especificar"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 800
$ws.Cells.Item(99, 11).Value = 2800
$ws.Cells.Item(99, 12).Value = 3000
$ws.Cells.Item(99, 13).Value = 2900
$ws.Cells.Item(99, 14).Value = "$/paquete"
$ws.Cells.Item(99, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(99, 16).Value = 2900
$ws.Cells.Item(99, 17).Value = 1
$ws.Cells.Item(99, 18).Value = "Hortaliza"
